# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 3 in the "Pepino dulce" sheet,
# shifting the existing rows 3-33 down to 4-34.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 3 (this shifts rows 3..33 -> 4..34)
$ws.Rows.Item(3).Insert()

# Populate the new row 3 with the new weekly record
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(3, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(3, 4).Value = 45043
$ws.Cells.Item(3, 5).Value = 15
$ws.Cells.Item(3, 6).Value = 100112043
$ws.Cells.Item(3, 7).Value = "Pepino dulce"
$ws.Cells.Item(3, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(3, 9).Value = "Segunda"
$ws.Cells.Item(3, 10).Value = 170
$ws.Cells.Item(3, 11).Value = 18000
$ws.Cells.Item(3, 12).Value = 20000
$ws.Cells.Item(3, 13).Value = 19059
$ws.Cells.Item(3, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(3, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(3, 16).Value = 1059
$ws.Cells.Item(3, 17).Value = 18
$ws.Cells.Item(3, 18).Value = "Hortaliza"
